# Add "NA" values under the duplicate_image_filename column (column E)
# for rows 2 through 21 (the rows that contain practice/test stimuli data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 21; $row++) {
    $ws.Range("E$row").Value = "NA"
}
